# Update the 4th "Featured News" headline (cell B5 on Sheet1) with the new
# news item text, replacing the old "Cognizant Flowsource..." headline with
# "Cognizant shines at the Times Group Global Business Summit (GBS)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Cognizant shines at the Times Group Global Business Summit (GBS)"
